$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value2 = "Volume 29   Number  43"
$ws.Range("C9").Value2 = "Report Covering the Week  10/24/2022  Through  10/30/2022"

# --- Type-changing cells: convert using copy/paste-special trick to pick up both value-type and style ---
# Source cells (row 30 is untouched by this edit, used purely as a style/type template):
#   C30 = text "0" w/ style 14 ; I30 = number w/ style 15 ; K30 = percent-number w/ style 16

# C15: number -> text "0" (style 14)
$ws.Range("C30").Copy()
$ws.Range("C15").PasteSpecial(-4123)
$ws.Range("C30").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# D15: text -> number 1
$ws.Range("I30").Copy()
$ws.Range("D15").PasteSpecial(-4123)
$ws.Range("I30").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value2 = 1

# E15: text -> number -100
$ws.Range("K30").Copy()
$ws.Range("E15").PasteSpecial(-4123)
$ws.Range("K30").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value2 = -100

# C20: number -> text "0" (style 14)
$ws.Range("C30").Copy()
$ws.Range("C20").PasteSpecial(-4123)
$ws.Range("C30").Copy()
$ws.Range("C20").PasteSpecial(-4122)

# D22: text -> number 1
$ws.Range("I30").Copy()
$ws.Range("D22").PasteSpecial(-4123)
$ws.Range("I30").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value2 = 1

# E22: text -> number -100
$ws.Range("K30").Copy()
$ws.Range("E22").PasteSpecial(-4123)
$ws.Range("K30").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value2 = -100

# C26: number -> text "0" (style 14)
$ws.Range("C30").Copy()
$ws.Range("C26").PasteSpecial(-4123)
$ws.Range("C30").Copy()
$ws.Range("C26").PasteSpecial(-4122)

# --- Simple same-type value updates ---
$ws.Range("I15").Value2 = 16
$ws.Range("J15").Value2 = 18
$ws.Range("K15").Value2 = -11.111111111111
$ws.Range("L15").Value2 = 14.285714285714
$ws.Range("M15").Value2 = 60
$ws.Range("N15").Value2 = -30.434782608695
$ws.Range("D16").Value2 = 6
$ws.Range("E16").Value2 = -66.666666666666
$ws.Range("G16").Value2 = 13
$ws.Range("H16").Value2 = -15.384615384615
$ws.Range("I16").Value2 = 98
$ws.Range("J16").Value2 = 90
$ws.Range("K16").Value2 = 8.888888888888
$ws.Range("L16").Value2 = 13.953488372093
$ws.Range("M16").Value2 = -31.944444444444
$ws.Range("N16").Value2 = -84.566929133858
$ws.Range("C17").Value2 = 3
$ws.Range("E17").Value2 = 50
$ws.Range("F17").Value2 = 19
$ws.Range("H17").Value2 = 72.727272727272
$ws.Range("I17").Value2 = 150
$ws.Range("J17").Value2 = 145
$ws.Range("K17").Value2 = 3.448275862068
$ws.Range("L17").Value2 = 27.118644067796
$ws.Range("M17").Value2 = 44.230769230769
$ws.Range("N17").Value2 = -44.029850746268
$ws.Range("C18").Value2 = 6
$ws.Range("D18").Value2 = 4
$ws.Range("E18").Value2 = 50
$ws.Range("F18").Value2 = 17
$ws.Range("G18").Value2 = 16
$ws.Range("H18").Value2 = 6.25
$ws.Range("I18").Value2 = 161
$ws.Range("J18").Value2 = 119
$ws.Range("K18").Value2 = 35.294117647058
$ws.Range("L18").Value2 = 35.294117647058
$ws.Range("M18").Value2 = -33.195020746888
$ws.Range("N18").Value2 = -88.047512991833
$ws.Range("C19").Value2 = 5
$ws.Range("D19").Value2 = 14
$ws.Range("E19").Value2 = -64.285714285714
$ws.Range("F19").Value2 = 51
$ws.Range("G19").Value2 = 49
$ws.Range("H19").Value2 = 4.081632653061
$ws.Range("I19").Value2 = 624
$ws.Range("J19").Value2 = 473
$ws.Range("K19").Value2 = 31.923890063424
$ws.Range("L19").Value2 = 61.658031088082
$ws.Range("M19").Value2 = 73.816155988857
$ws.Range("N19").Value2 = -4
$ws.Range("D20").Value2 = 3
$ws.Range("E20").Value2 = -100
$ws.Range("F20").Value2 = 10
$ws.Range("G20").Value2 = 18
$ws.Range("H20").Value2 = -44.444444444444
$ws.Range("J20").Value2 = 91
$ws.Range("K20").Value2 = 37.362637362637
$ws.Range("L20").Value2 = 52.439024390243
$ws.Range("M20").Value2 = -6.716417910447
$ws.Range("N20").Value2 = -91.694352159468
$ws.Range("C21").Value2 = 16
$ws.Range("D21").Value2 = 30
$ws.Range("E21").Value2 = -46.666666666666
$ws.Range("F21").Value2 = 110
$ws.Range("G21").Value2 = 110
$ws.Range("H21").Value2 = 0
$ws.Range("I21").Value2 = 1179
$ws.Range("J21").Value2 = 938
$ws.Range("K21").Value2 = 25.692963752665
$ws.Range("L21").Value2 = 46.096654275092
$ws.Range("M21").Value2 = 18.136272545090
$ws.Range("N21").Value2 = -73.427991886409
$ws.Range("G22").Value2 = 2
$ws.Range("H22").Value2 = -50
$ws.Range("J22").Value2 = 11
$ws.Range("K22").Value2 = 0
$ws.Range("M22").Value2 = -54.166666666666
$ws.Range("C24").Value2 = 37
$ws.Range("D24").Value2 = 27
$ws.Range("E24").Value2 = 37.037037037037
$ws.Range("F24").Value2 = 151
$ws.Range("H24").Value2 = 57.291666666666
$ws.Range("I24").Value2 = 1553
$ws.Range("J24").Value2 = 1005
$ws.Range("K24").Value2 = 54.527363184079
$ws.Range("L24").Value2 = 52.854330708661
$ws.Range("M24").Value2 = 86.658653846153
$ws.Range("C25").Value2 = 8
$ws.Range("D25").Value2 = 11
$ws.Range("E25").Value2 = -27.272727272727
$ws.Range("F25").Value2 = 40
$ws.Range("H25").Value2 = 17.647058823529
$ws.Range("I25").Value2 = 380
$ws.Range("J25").Value2 = 360
$ws.Range("K25").Value2 = 5.555555555555
$ws.Range("L25").Value2 = 34.275618374558
$ws.Range("M25").Value2 = 2.425876010781
$ws.Range("D26").Value2 = 1
$ws.Range("E26").Value2 = -100
$ws.Range("I26").Value2 = 21
$ws.Range("J26").Value2 = 22
$ws.Range("K26").Value2 = -4.545454545454
$ws.Range("L26").Value2 = -12.5
$ws.Range("C27").Value2 = 3
$ws.Range("D27").Value2 = 2
$ws.Range("E27").Value2 = 50
$ws.Range("F27").Value2 = 7
$ws.Range("G27").Value2 = 9
$ws.Range("H27").Value2 = -22.222222222222
$ws.Range("I27").Value2 = 46
$ws.Range("J27").Value2 = 52
$ws.Range("K27").Value2 = -11.538461538461
$ws.Range("L27").Value2 = 53.333333333333
$ws.Range("L28").Value2 = 66.666666666666
$ws.Range("L29").Value2 = 66.666666666666
